$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Layout change -------------------------------------------------------
# A new station "S8" (and its shadow column "eS8") was added to the report.
# Insert a blank column before the old "Total" column (col 10 / J) for the
# new "S8" series, and another blank column before the old "eTotal" column
# (which, after the first insert, sits at col 19 / S) for "eS8".
$ws.Columns.Item(10).Insert()
$ws.Columns.Item(19).Insert()

# Label the two newly-inserted header cells (they inherit the bold header
# style automatically from the Insert, matching the other header cells).
$ws.Cells.Item(1, 10).Value = "S8"
$ws.Cells.Item(1, 19).Value = "eS8"

# --- Data ------------------------------------------------------------------
# Columns: A=Date B=Hour C=S1 D=S2 E=S3 F=S4 G=S5 H=S6 I=S7 J=S8 K=Total
#          L=eS1 M=eS2 N=eS3 O=eS4 P=eS5 Q=eS6 R=eS7 S=eS8 T=eTotal
$data = @(
    @("2024-11-07", 1, 20222, 5513, 11536, 16729, 0, 0, 0, 0, 54000, 20224.0222, 5513.5513, 11537.1536, 16730.6729, 0, 0, 0, 0, 54005.4),
    @("2024-11-07", 2, 19240, 5184, 11048, 15329, 0, 0, 0, 0, 50801, 19241.924, 5184.5184, 11049.1048, 15330.5329, 0, 0, 0, 0, 50806.0801),
    @("2024-11-07", 3, 18198, 4854, 10610, 14949, 0, 0, 0, 0, 48611, 18199.8198, 4854.4854, 10611.061, 14950.4949, 0, 0, 0, 0, 48615.8611),
    @("2024-11-07", 4, 17345, 4658, 10252, 14469, 0, 0, 0, 0, 46724, 17346.7345, 4658.4658, 10253.0252, 14470.4469, 0, 0, 0, 0, 46728.6724),
    @("2024-11-07", 5, 17742, 4800, 10535, 14489, 0, 0, 0, 0, 47566, 17743.7742, 4800.48, 10536.0535, 14490.4489, 0, 0, 0, 0, 47570.75659999999),
    @("2024-11-07", 6, 18085, 4890, 10899, 14611, 0, 0, 0, 0, 48485, 18086.8085, 4890.489, 10900.0899, 14612.4611, 0, 0, 0, 0, 48489.8485),
    @("2024-11-07", 7, 17689, 4420, 10825, 14937, 0, 0, 0, 0, 47871, 17690.7689, 4420.442, 10826.0825, 14938.4937, 0, 0, 0, 0, 47875.7871),
    @("2024-11-07", 8, 20608, 4502, 0, 16873, 0, 0, 0, 0, 41983, 20610.0608, 4502.4502, 0, 16874.6873, 0, 0, 0, 0, 41987.1983),
    @("2024-11-07", 9, 25888, 5559, 14619, 19334, 0, 0, 0, 0, 65400, 25890.5888, 5559.5559, 14620.4619, 19335.9334, 0, 0, 0, 0, 65406.53999999999),
    @("2024-11-07", 10, 28174, 5754, 19760, 20151, 0, 0, 0, 0, 73839, 28176.8174, 5754.5754, 19761.976, 20153.0151, 0, 0, 0, 0, 73846.3839),
    @("2024-11-07", 11, 29407, 6128, 21229, 20872, 0, 0, 0, 0, 77636, 29409.9407, 6128.6128, 21231.1229, 20874.0872, 0, 0, 0, 0, 77643.76359999999),
    @("2024-11-07", 12, 30004, 6245, 20824, 21493, 0, 0, 0, 0, 78566, 30007.0004, 6245.6245, 20826.0824, 21495.1493, 0, 0, 0, 0, 78573.8566),
    @("2024-11-07", 13, 30232, 6142, 21177, 21767, 0, 0, 0, 0, 79318, 30235.0232, 6142.6142, 21179.1177, 21769.1767, 0, 0, 0, 0, 79325.93179999999),
    @("2024-11-07", 14, 31909, 6459, 21970, 22757, 0, 0, 0, 0, 83095, 31912.1909, 6459.6459, 21972.197, 22759.2757, 0, 0, 0, 0, 83103.3095),
    @("2024-11-07", 15, 0, 6674, 22172, 23064, 0, 0, 0, 0, 51910, 0, 6674.6674, 22174.2172, 23066.3064, 0, 0, 0, 0, 51915.191),
)

# The "Date" column holds plain text ("2024-11-07"), not a real date. Force
# text formatting before assigning so Excel doesn't auto-convert the string
# into a date serial value, then clear the temporary format again so the
# cell ends up with the same (default) style as the rest of the data rows.
$dateRng = $ws.Range("A2:A16")
$dateRng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

$dateRng.ClearFormats()
